# "Test: Add prime attribute test"
# Adds a new "Prime" worksheet (after "AttributeLimit") containing a
# black-box test case table for the "max strength attribute" / prime
# attribute rule, mirroring the layout used on the other test sheets.

$wb = $excel.ActiveWorkbook

# --- restore the per-sheet selections the author left on the existing
#     sheets (each worksheet remembers its own last selection; whichever
#     sheet we touch last becomes the active tab) --------------------
$wsConstructor = $wb.Worksheets.Item("Constructor")
$wsConstructor.Range("D4").Select()

$wsCombatMastery = $wb.Worksheets.Item("CombatMastery")
$wsCombatMastery.Range("C7").Select()

$wsAttributeLimit = $wb.Worksheets.Item("AttributeLimit")
$wsAttributeLimit.Range("A1:F1").Select()

# --- add the new "Prime" sheet after the last existing sheet --------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Prime"

# --- header rows (1 & 2): bold, centered both ways -------------------
$headerRange = $ws.Range("A1:J2")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4108     # xlCenter

$ws.Range("A1").Value = "Test Case"
$ws.Range("B1").Value = "Category"
$ws.Range("C1").Value = "Partition"
$ws.Range("D1").Value = "Test Inputs"
$ws.Range("I1").Value = "Expected Output"
$ws.Range("J1").Value = "Comments"

# --- data row 3 values first, so shared-string insertion order matches
#     ("Max strength attribute" needs to land before the attribute
#     sub-headers below) ------------------------------------------
$dataRange = $ws.Range("A3:I3")
$dataRange.HorizontalAlignment = -4108     # xlCenter, no vertical centering

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Valid Input"
$ws.Range("C3").Value = "Max strength attribute"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 3

# --- attribute sub-headers on row 2 (italic, centered both ways) ----
$subHeaderRange = $ws.Range("D2:H2")
$subHeaderRange.Font.Italic = $true
$subHeaderRange.HorizontalAlignment = -4108
$subHeaderRange.VerticalAlignment = -4108

$ws.Range("D2").Value = "Level"
$ws.Range("E2").Value = "Might"
$ws.Range("F2").Value = "Agility"
$ws.Range("G2").Value = "Intelligence"
$ws.Range("H2").Value = "Charisma"

# --- merge the header cells that span both rows / multiple columns --
$ws.Range("D1:H1").Merge()
$ws.Range("A1:A2").Merge()
$ws.Range("B1:B2").Merge()
$ws.Range("C1:C2").Merge()
$ws.Range("I1:I2").Merge()
$ws.Range("J1:J2").Merge()

# --- leave the new sheet's selection/active-cell where the author's
#     saved file shows it -------------------------------------------
$ws.Range("J3").Select()
